$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WDC")

# Row 4 - Inventory
$ws.Range("C4").Value = 3576000000.0
$ws.Range("D4").Value = 3355000000.0
$ws.Range("E4").Value = 3070000000.0
$ws.Range("F4").Value = 3091000000.0
$ws.Range("G4").Value = 3122000000.0

# Row 14 - Accounts Payable
$ws.Range("C14").Value = 2332000000.0
$ws.Range("D14").Value = 2353000000.0
$ws.Range("E14").Value = 2352000000.0
$ws.Range("F14").Value = 2183000000.0
$ws.Range("G14").Value = 2100000000.0

# Row 20 - Long Term Tax Liability (Deferred)
$ws.Range("C20").Value = 724000000.0
$ws.Range("D20").Value = 715000000.0
$ws.Range("E20").Value = 720000000.0
$ws.Range("F20").Value = 716000000.0
$ws.Range("G20").Value = 708000000.0

# Row 33 - Net Debt
$ws.Range("B33").Value = 6195000000.0

# Row 34 - Total Debt
$ws.Range("B34").Value = 8929000000.0
